$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.118.83"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.107.91"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "348.89"
$ws.Range("E5").Value = "  +3.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4446"
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.69"
$ws.Range("E9").Value = "  -4.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08965"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.174"
$ws.Range("E11").Value = "  +0.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.71"
$ws.Range("E12").Value = "  +4.65%  "
$ws.Range("D13").Value = "2.111.02"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.228"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.748"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.23"
$ws.Range("E16").Value = "  +2.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001149"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.84"
$ws.Range("E19").Value = "  +7.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06696"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.232"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "30.220.42"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.85"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.341"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "2.359.24"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.00"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.538"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.99"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.73"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.176"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1068"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.634"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.255"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.975"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.44"
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.935"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02577"
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06828"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2307"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.59"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6828"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.278"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.24"
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.309"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6377"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000367"
$ws.Range("E47").Value = "  +3.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.656"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.223"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.36"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07240"
$ws.Range("E51").Value = "  +0.51%  "
